# Update column G ("K") values on Sheet1 with newly regenerated figures,
# replacing the old "Strike#"-derived numbers (commit: "regen save_data to
# use K instead of Strike#, regen std/mean, calc and write s_vals").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 0
    3  = 0
    4  = 1
    5  = 0
    6  = 2
    7  = 0
    8  = 1
    9  = 1
    10 = 0
    11 = 2
    12 = 2
    13 = 3
    14 = 0
    15 = 3
    16 = 1
    17 = 2
    18 = 2
    19 = 3
    20 = 0
    21 = 1
    22 = 1
    23 = 2
    24 = 1
    25 = 1
    26 = 3
    27 = 1
    28 = 1
    29 = 0
    30 = 1
    31 = 0
    32 = 3
    33 = 0
    34 = 0
    35 = 3
    36 = 4
    37 = 2
    38 = 1
    39 = 1
    40 = 0
    41 = 2
    42 = 2
    43 = 0
    44 = 0
    45 = 2
    46 = 0
    47 = 1
    48 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
